$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Re-colour existing rows 19, 51 and 52 (requirement rows that were
#    reclassified) by copying the *format* from rows that already carry the
#    target cell style, so no new style records are created.
# ---------------------------------------------------------------------------
$ws.Range("B17").Copy() | Out-Null
$ws.Range("A19:E19").PasteSpecial(-4122) | Out-Null

$ws.Range("E5").Copy() | Out-Null
$ws.Range("A51:E52").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Append six new "Search" requirement rows (171-176) at the bottom of the
#    sheet, two of which (B176/D176) store their numbers as text - matching
#    the existing REQ_ID / DEPENDENCY column conventions.
# ---------------------------------------------------------------------------

# Columns A, B and D carry a column-level Text ("@") number format, so a
# plain `.Value = <number>` assignment on them would be stored as a string.
# Write genuine numbers into a scratch cell outside the formatted columns
# first, then *move* (Cut/Paste) them into place - a move preserves the
# literal numeric type instead of re-parsing it under the destination
# format.
function Set-NumericValue($cellAddress, $number) {
    $scratch = $ws.Range("Z1")
    $scratch.Value = $number
    $scratch.Cut($ws.Range($cellAddress)) | Out-Null
}

Set-NumericValue "B171" 155
$ws.Range("C171").Value = "Search"
Set-NumericValue "D171" 4
$ws.Range("E171").Value = "Searching for an empty string returns a comprehensive list of all movie and people results"

Set-NumericValue "B172" 156
$ws.Range("C172").Value = "Search"
Set-NumericValue "D172" 4
$ws.Range("E172").Value = "When searching for a string, if any movie matches are found, then only the movie results are displayed."

Set-NumericValue "B173" 157
$ws.Range("C173").Value = "Search"
Set-NumericValue "D173" 4
$ws.Range("E173").Value = "When searching for a string, if no movie matches are found then any people matches are shown instead."

Set-NumericValue "B174" 158
$ws.Range("C174").Value = "Search"
Set-NumericValue "D174" 4
$ws.Range("E174").Value = "When searching for a string, if no movie or people matches are found, then display a comprehensive list of all movie and people results."

Set-NumericValue "B175" 159
$ws.Range("C175").Value = "Search"
Set-NumericValue "D175" 4
$ws.Range("E175").Value = "When searching for a string, a substring can be a match."

# Row 176 - B/D are textual ("160" / "4"), matching other rows (e.g. 162/163)
# that store REQ_ID numbers as text once double-digit/irregular formatting
# kicks in.
$ws.Range("E176").Value = "When searching for a string, if the string contains more than one word then individual words can be a match."
$ws.Range("B176").Value = "160"
$ws.Range("D176").Value = "4"
$ws.Range("C176").Value = "Search"

# ---------------------------------------------------------------------------
# 3) Apply the correct cell styles to the new rows by copying formats from
#    cells that already use the desired style index (keeps styles.xml
#    untouched instead of registering new xf records).
# ---------------------------------------------------------------------------
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B171:B175").PasteSpecial(-4122) | Out-Null
$ws.Range("D171:D175").PasteSpecial(-4122) | Out-Null

$ws.Range("C17").Copy() | Out-Null
$ws.Range("C171:C176").PasteSpecial(-4122) | Out-Null
$ws.Range("E171:E175").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Copy() | Out-Null
$ws.Range("B176:D176").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Clear the scratch cell used to stage numeric literals.
$ws.Range("Z1").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 4) Update the sheet view to focus on the newly-added rows.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 150
$ws.Range("D171:D175").Select()
